$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2000
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1830

$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -16

$ws.Range("H32").Value = 10115.375
$ws.Range("I32").Value = 12429
$ws.Range("J32").Value = 9162.706
$ws.Range("K32").Value = 12429
$ws.Range("L32").Value = 9162.706
$ws.Range("M32").Value = -12103
$ws.Range("N32").Value = -9814.706

$ws.Range("H53").Value = 356
$ws.Range("I53").Value = 250
$ws.Range("J53").Value = 462
$ws.Range("K53").Value = 250
$ws.Range("L53").Value = 462
$ws.Range("M53").Value = 387
$ws.Range("N53").Value = -1736

$ws.Range("H64").Value = 11355.333
$ws.Range("I64").Value = 8799.5
$ws.Range("J64").Value = 13400
$ws.Range("K64").Value = 8799.5
$ws.Range("L64").Value = 13400
$ws.Range("M64").Value = -8551.5
$ws.Range("N64").Value = -13896

$ws.Range("H67").Value = 11355.333
$ws.Range("I67").Value = 8799.5
$ws.Range("J67").Value = 13400
$ws.Range("K67").Value = 8799.5
$ws.Range("L67").Value = 13400
$ws.Range("M67").Value = -7941.5
$ws.Range("N67").Value = -15116

$ws.Range("H69").Value = 7975.143
$ws.Range("I69").Value = 1997.5
$ws.Range("J69").Value = 10366.2
$ws.Range("K69").Value = 5992.5
$ws.Range("L69").Value = 31098.6
$ws.Range("M69").Value = -5118.5
$ws.Range("N69").Value = -32846.60000000001

$ws.Range("H72").Value = 7975.143
$ws.Range("I72").Value = 1997.5
$ws.Range("J72").Value = 10366.2
$ws.Range("K72").Value = 17977.5
$ws.Range("L72").Value = 93295.8
$ws.Range("M72").Value = -13609.5
$ws.Range("N72").Value = -102031.8

$ws.Range("H76").Value = 4839
$ws.Range("I76").Value = 5499
$ws.Range("J76").Value = 4729
$ws.Range("K76").Value = 5499
$ws.Range("L76").Value = 4729
$ws.Range("M76").Value = -5184
$ws.Range("N76").Value = -5359

$ws.Range("H79").Value = 4839
$ws.Range("I79").Value = 5499
$ws.Range("J79").Value = 4729
$ws.Range("K79").Value = 5499
$ws.Range("L79").Value = 4729
$ws.Range("M79").Value = -4407
$ws.Range("N79").Value = -6913

$ws.Range("H98").Value = 951.8570999999999
$ws.Range("I98").Value = 903.6667
$ws.Range("J98").Value = 1038.6
$ws.Range("K98").Value = 903.6667
$ws.Range("L98").Value = 1038.6
$ws.Range("M98").Value = 594.3333
$ws.Range("N98").Value = -4034.6

$ws.Range("H99").Value = 2859
$ws.Range("I99").Value = 424.57144
$ws.Range("J99").Value = 8539.333000000001
$ws.Range("K99").Value = 1273.71432
$ws.Range("L99").Value = 25617.999
$ws.Range("M99").Value = 224.28568
$ws.Range("N99").Value = -28613.999

$ws.Range("H106").Value = 1074.7858
$ws.Range("I106").Value = 857.46155
$ws.Range("J106").Value = 3900
$ws.Range("K106").Value = 857.46155
$ws.Range("L106").Value = 3900
$ws.Range("M106").Value = -226.46155
$ws.Range("N106").Value = -5162

$ws.Range("H107").Value = 560.13794
$ws.Range("I107").Value = 573.5925999999999
$ws.Range("J107").Value = 378.5
$ws.Range("K107").Value = 573.5925999999999
$ws.Range("L107").Value = 378.5
$ws.Range("M107").Value = 1346.4074
$ws.Range("N107").Value = -4218.5

$ws.Range("H112").Value = 101669.5
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 112843.89
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 338531.67
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -340747.67

$ws.Range("H122").Value = 951.8570999999999
$ws.Range("I122").Value = 903.6667
$ws.Range("J122").Value = 1038.6
$ws.Range("K122").Value = 2711.0001
$ws.Range("L122").Value = 3115.8
$ws.Range("M122").Value = -261.0001000000002
$ws.Range("N122").Value = -8015.799999999999

$ws.Range("H137").Value = 1364.7097
$ws.Range("I137").Value = 1314.9259
$ws.Range("J137").Value = 1700.75
$ws.Range("K137").Value = 3944.7777
$ws.Range("L137").Value = 5102.25
$ws.Range("M137").Value = -1394.7777
$ws.Range("N137").Value = -10202.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3074.9412
$ws.Range("I2").Value = 2767.125
$ws.Range("J2").Value = 8000
$ws.Range("K2").Value = 2767.125
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = -2654.125
$ws.Range("N2").Value = -8226

$ws.Range("H3").Value = 67000
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 100000
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 100000
$ws.Range("M3").Value = -885
$ws.Range("N3").Value = -100230

$ws.Range("H36").Value = 44726.43
$ws.Range("I36").Value = 23017.2
$ws.Range("J36").Value = 98999.5
$ws.Range("K36").Value = 23017.2
$ws.Range("L36").Value = 98999.5
$ws.Range("M36").Value = -22671.2
$ws.Range("N36").Value = -99691.5

$ws.Range("H97").Value = 943.25
$ws.Range("I97").Value = 1123.88
$ws.Range("J97").Value = 298.14285
$ws.Range("K97").Value = 1123.88
$ws.Range("L97").Value = 298.14285
$ws.Range("M97").Value = -627.8800000000001
$ws.Range("N97").Value = -1290.14285

$ws.Range("H104").Value = 59000
$ws.Range("I104").Value = 1000
$ws.Range("J104").Value = 117000
$ws.Range("K104").Value = 1000
$ws.Range("L104").Value = 117000
$ws.Range("M104").Value = 2494
$ws.Range("N104").Value = -123988

$ws.Range("H116").Value = 3074.9412
$ws.Range("I116").Value = 2767.125
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 2767.125
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -473.125
$ws.Range("N116").Value = -12588

$ws.Range("H135").Value = 79950
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 79950
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 79950
$ws.Range("N135").Value = -90090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3074.9412
$ws.Range("I3").Value = 2767.125
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 2767.125
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -2653.125
$ws.Range("N3").Value = -8228

$ws.Range("H94").Value = 1838.4814
$ws.Range("I94").Value = 1539.5416
$ws.Range("J94").Value = 4230
$ws.Range("K94").Value = 1539.5416
$ws.Range("L94").Value = 4230
$ws.Range("M94").Value = -1088.5416
$ws.Range("N94").Value = -5132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 900
$ws.Range("I16").Value = 900
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -613
$ws.Range("N16").ClearContents()

$ws.Range("H55").Value = 10710.286
$ws.Range("I55").Value = 9162
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 9162
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -8847
$ws.Range("N55").Value = -20630

$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1270
$ws.Range("N113").ClearContents()

$ws.Range("H135").Value = 74315
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 74315
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 74315
$ws.Range("N135").Value = -84455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 153.45454
$ws.Range("I38").Value = 171.42857
$ws.Range("J38").Value = 122
$ws.Range("K38").Value = 514.28571
$ws.Range("L38").Value = 366
$ws.Range("M38").Value = -167.28571
$ws.Range("N38").Value = -1060

$ws.Range("H40").Value = 312.8095
$ws.Range("I40").Value = 198.27272
$ws.Range("J40").Value = 438.8
$ws.Range("K40").Value = 793.09088
$ws.Range("L40").Value = 1755.2
$ws.Range("M40").Value = -724.09088
$ws.Range("N40").Value = -1893.2

$ws.Range("H86").Value = 607.8
$ws.Range("I86").Value = 566.6667
$ws.Range("J86").Value = 669.5
$ws.Range("K86").Value = 1700.0001
$ws.Range("L86").Value = 2008.5
$ws.Range("M86").Value = -514.0001
$ws.Range("N86").Value = -4380.5

$ws.Range("H89").Value = 607.8
$ws.Range("I89").Value = 566.6667
$ws.Range("J89").Value = 669.5
$ws.Range("K89").Value = 5100.0003
$ws.Range("L89").Value = 6025.5
$ws.Range("M89").Value = 827.9997000000003
$ws.Range("N89").Value = -17881.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 200419.28
$ws.Range("I42").Value = 200000
$ws.Range("J42").Value = 200587
$ws.Range("K42").Value = 200000
$ws.Range("L42").Value = 200587
$ws.Range("M42").Value = -199515
$ws.Range("N42").Value = -201557

$ws.Range("H105").Value = 26789
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 26789
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 26789
$ws.Range("N105").Value = -33777

$ws.Range("H115").Value = 200419.28
$ws.Range("I115").Value = 200000
$ws.Range("J115").Value = 200587
$ws.Range("K115").Value = 200000
$ws.Range("L115").Value = 200587
$ws.Range("M115").Value = -198825
$ws.Range("N115").Value = -202937

$ws.Range("H122").Value = 2124.125
$ws.Range("I122").Value = 2089.5833
$ws.Range("J122").Value = 2227.75
$ws.Range("K122").Value = 6268.749899999999
$ws.Range("L122").Value = 6683.25
$ws.Range("M122").Value = -3818.749899999999
$ws.Range("N122").Value = -11583.25

$ws.Range("H126").Value = 17360.045
$ws.Range("I126").Value = 27017.154
$ws.Range("J126").Value = 3410.889
$ws.Range("K126").Value = 81051.462
$ws.Range("L126").Value = 10232.667
$ws.Range("M126").Value = -78581.462
$ws.Range("N126").Value = -15172.667

$ws.Range("H135").Value = 40546
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 40546
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40546
$ws.Range("N135").Value = -50686

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1142.3077
$ws.Range("I22").Value = 890.625
$ws.Range("J22").Value = 1545
$ws.Range("K22").Value = 890.625
$ws.Range("L22").Value = 1545
$ws.Range("M22").Value = -595.625
$ws.Range("N22").Value = -2135

$ws.Range("H27").Value = 1142.3077
$ws.Range("I27").Value = 890.625
$ws.Range("J27").Value = 1545
$ws.Range("K27").Value = 890.625
$ws.Range("L27").Value = 1545
$ws.Range("M27").Value = -783.625
$ws.Range("N27").Value = -1759

$ws.Range("H55").Value = 358.1111
$ws.Range("I55").Value = 422.75
$ws.Range("J55").Value = 306.4
$ws.Range("K55").Value = 422.75
$ws.Range("L55").Value = 306.4
$ws.Range("M55").Value = -249.75
$ws.Range("N55").Value = -652.4

$ws.Range("H128").Value = 25000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 25000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 25000
$ws.Range("N128").Value = -34960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 60500
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 60500
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 60500
$ws.Range("N16").Value = -61084

$ws.Range("H62").Value = 18033.812
$ws.Range("I62").Value = 22354.2
$ws.Range("J62").Value = 10833.167
$ws.Range("K62").Value = 22354.2
$ws.Range("L62").Value = 10833.167
$ws.Range("M62").Value = -21730.2
$ws.Range("N62").Value = -12081.167

$ws.Range("H65").Value = 18033.812
$ws.Range("I65").Value = 22354.2
$ws.Range("J65").Value = 10833.167
$ws.Range("K65").Value = 111771
$ws.Range("L65").Value = 54165.835
$ws.Range("M65").Value = -108651
$ws.Range("N65").Value = -60405.835

$ws.Range("H122").Value = 2758.3914
$ws.Range("I122").Value = 2707.9
$ws.Range("J122").Value = 2853.0625
$ws.Range("K122").Value = 8123.700000000001
$ws.Range("L122").Value = 8559.1875
$ws.Range("M122").Value = -5673.700000000001
$ws.Range("N122").Value = -13459.1875

$ws.Range("H136").Value = 2028.4166
$ws.Range("I136").Value = 1938
$ws.Range("J136").Value = 2589
$ws.Range("K136").Value = 5814
$ws.Range("L136").Value = 7767
$ws.Range("M136").Value = -3264
$ws.Range("N136").Value = -12867
